# Generate Report for Archive
$wb = $excel.ActiveWorkbook

# Update the status text across all sheets that use the shared string
# "Ready for handoff" -> "In Translation"
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $val = $cell.Value2
        if (($val -is [string]) -and ($val -eq "Ready for handoff")) {
            $cell.Value = "In Translation"
        }
    }
}

# Narrow the columns that held the status text so their width matches the
# shorter string ("In Translation" is shorter than "Ready for handoff").
# (Target stored column width is ~13.41; ColumnWidth is quantized to the
# pixel grid on write, so 12.5 is the input that lands closest to it.)
$newWidth = 12.5

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = $newWidth

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = $newWidth
